$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Years"
$ws.Range("G2").Value = "1958, 1962, 1970, 1994, 2002"
$ws.Range("G3").Value = "1954, 1974, 1990, 2014"
$ws.Range("G4").Value = "1934, 1938, 1982,  2006"
$ws.Range("G5").Value = "1978, 1986"
$ws.Range("G6").Value = "1930, 1950"
$ws.Range("G7").Value = 1998
$ws.Range("G8").Value = 1966
$ws.Range("G9").Value = 2010

$ws.Range("G10").Select()
